# Update workbook with new daily rows (update through 26/04/2021)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Serial date values for the 5 new rows (continuing the daily series)
$dates = @(44308, 44309, 44310, 44311, 44312)

$startRow = 234
$row = $startRow
foreach ($d in $dates) {
    $ws.Cells.Item($row, 1).Value = $d
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
    $row = $row + 1
}

# Copy the formatting (style) of the last existing data row (A233) onto the
# date column of the new rows so they keep the same date style (s="2").
$ws.Range("A233").Copy()
$ws.Range("A234:A238").PasteSpecial(-4122)
